$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.232.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "'1.604.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'18.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "'1.827.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "'1.606.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "'26.205.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "'62.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D20").Value = "'200.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'9.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'5.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "'144.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("D28").Value = "'15.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").Value = "'0.0490"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.87%  "
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "'1.163.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'2.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "'0.784"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "'0.497"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "'5.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").Value = "'1.739.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'91.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("E46").Value = "  +15.47%  "
$ws.Range("D47").Value = "'1.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  -0.03%  "
